$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix P20: was a suppressed "<3" marker, now a hard 0 ---
$ws.Range("P20").Value = 0

# --- Fill in full data for rows 21-35 (MARYLAND .. NORTH DAKOTA) ---
    # Row 21
    $ws.Range("B21").Value = 93
    $ws.Range("C21").Value = 17
    $ws.Range("D21").Value = 1
    $ws.Range("E21").Value = 7
    $ws.Range("F21").Value = 118
    $ws.Range("G21").Value = 15
    $ws.Range("H21").Value = "<3"
    $ws.Range("I21").Value = 11
    $ws.Range("J21").Value = 5
    $ws.Range("K21").Value = 0
    $ws.Range("L21").Value = 59
    $ws.Range("M21").Value = 29
    $ws.Range("N21").Value = 87
    $ws.Range("O21").Value = "<3"
    $ws.Range("P21").Value = 0
    # Row 22
    $ws.Range("B22").Value = 51
    $ws.Range("C22").Value = 4
    $ws.Range("D22").Value = 1
    $ws.Range("E22").Value = 6
    $ws.Range("F22").Value = 62
    $ws.Range("G22").Value = 4
    $ws.Range("H22").Value = "<3"
    $ws.Range("I22").Value = 3
    $ws.Range("J22").Value = "<3"
    $ws.Range("K22").Value = 0
    $ws.Range("L22").Value = 36
    $ws.Range("M22").Value = 17
    $ws.Range("N22").Value = 45
    $ws.Range("O22").Value = 5
    $ws.Range("P22").Value = 0
    # Row 23
    $ws.Range("B23").Value = 122
    $ws.Range("C23").Value = 18
    $ws.Range("D23").Value = 7
    $ws.Range("E23").Value = 5
    $ws.Range("F23").Value = 152
    $ws.Range("G23").Value = 16
    $ws.Range("H23").Value = "<3"
    $ws.Range("I23").Value = 9
    $ws.Range("J23").Value = 10
    $ws.Range("K23").Value = 0
    $ws.Range("L23").Value = 66
    $ws.Range("M23").Value = 48
    $ws.Range("N23").Value = 107
    $ws.Range("O23").Value = 6
    $ws.Range("P23").Value = 0
    # Row 24
    $ws.Range("B24").Value = 27
    $ws.Range("C24").Value = 8
    $ws.Range("D24").Value = 0
    $ws.Range("E24").Value = 1
    $ws.Range("F24").Value = 36
    $ws.Range("G24").Value = 8
    $ws.Range("H24").Value = 0
    $ws.Range("I24").Value = 5
    $ws.Range("J24").Value = 4
    $ws.Range("K24").Value = 0
    $ws.Range("L24").Value = 18
    $ws.Range("M24").Value = 10
    $ws.Range("N24").Value = 25
    $ws.Range("O24").Value = "<3"
    $ws.Range("P24").Value = 0
    # Row 25
    $ws.Range("B25").Value = 26
    $ws.Range("C25").Value = 5
    $ws.Range("D25").Value = 2
    $ws.Range("E25").Value = 3
    $ws.Range("F25").Value = 36
    $ws.Range("G25").Value = 4
    $ws.Range("H25").Value = "<3"
    $ws.Range("I25").Value = "<3"
    $ws.Range("J25").Value = 3
    $ws.Range("K25").Value = 0
    $ws.Range("L25").Value = 18
    $ws.Range("M25").Value = 7
    $ws.Range("N25").Value = 24
    $ws.Range("O25").Value = "<3"
    $ws.Range("P25").Value = 0
    # Row 26
    $ws.Range("B26").Value = 58
    $ws.Range("C26").Value = 4
    $ws.Range("D26").Value = 1
    $ws.Range("E26").Value = 4
    $ws.Range("F26").Value = 67
    $ws.Range("G26").Value = 4
    $ws.Range("H26").Value = "<3"
    $ws.Range("I26").Value = 4
    $ws.Range("J26").Value = "<3"
    $ws.Range("K26").Value = 0
    $ws.Range("L26").Value = 36
    $ws.Range("M26").Value = 20
    $ws.Range("N26").Value = 54
    $ws.Range("O26").Value = 3
    $ws.Range("P26").Value = 0
    # Row 27
    $ws.Range("B27").Value = 11
    $ws.Range("C27").Value = 5
    $ws.Range("D27").Value = 0
    $ws.Range("E27").Value = 3
    $ws.Range("F27").Value = 19
    $ws.Range("G27").Value = 5
    $ws.Range("H27").Value = 0
    $ws.Range("I27").Value = "<3"
    $ws.Range("J27").Value = 3
    $ws.Range("K27").Value = 0
    $ws.Range("L27").Value = 10
    $ws.Range("M27").Value = "<3"
    $ws.Range("N27").Value = 11
    $ws.Range("O27").Value = "<3"
    $ws.Range("P27").Value = 0
    # Row 28
    $ws.Range("B28").Value = 14
    $ws.Range("C28").Value = 7
    $ws.Range("D28").Value = 1
    $ws.Range("E28").Value = 2
    $ws.Range("F28").Value = 24
    $ws.Range("G28").Value = 6
    $ws.Range("H28").Value = "<3"
    $ws.Range("I28").Value = "<3"
    $ws.Range("J28").Value = 5
    $ws.Range("K28").Value = 0
    $ws.Range("L28").Value = 7
    $ws.Range("M28").Value = 6
    $ws.Range("N28").Value = 14
    $ws.Range("O28").Value = 0
    $ws.Range("P28").Value = 0
    # Row 29
    $ws.Range("B29").Value = 114
    $ws.Range("C29").Value = 11
    $ws.Range("D29").Value = 4
    $ws.Range("E29").Value = 4
    $ws.Range("F29").Value = 133
    $ws.Range("G29").Value = 10
    $ws.Range("H29").Value = "<3"
    $ws.Range("I29").Value = 5
    $ws.Range("J29").Value = 7
    $ws.Range("K29").Value = 0
    $ws.Range("L29").Value = 85
    $ws.Range("M29").Value = 31
    $ws.Range("N29").Value = 107
    $ws.Range("O29").Value = 6
    $ws.Range("P29").Value = 0
    # Row 30
    $ws.Range("B30").Value = 8
    $ws.Range("C30").Value = 3
    $ws.Range("D30").Value = 0
    $ws.Range("E30").Value = 1
    $ws.Range("F30").Value = 12
    $ws.Range("G30").Value = 3
    $ws.Range("H30").Value = 0
    $ws.Range("I30").Value = "<3"
    $ws.Range("J30").Value = "<3"
    $ws.Range("K30").Value = 0
    $ws.Range("L30").Value = 7
    $ws.Range("M30").Value = "<3"
    $ws.Range("N30").Value = 7
    $ws.Range("O30").Value = "<3"
    $ws.Range("P30").Value = 0
    # Row 31
    $ws.Range("B31").Value = 139
    $ws.Range("C31").Value = 13
    $ws.Range("D31").Value = 4
    $ws.Range("E31").Value = 13
    $ws.Range("F31").Value = 169
    $ws.Range("G31").Value = 13
    $ws.Range("H31").Value = "<3"
    $ws.Range("I31").Value = 8
    $ws.Range("J31").Value = 6
    $ws.Range("K31").Value = 0
    $ws.Range("L31").Value = 57
    $ws.Range("M31").Value = 75
    $ws.Range("N31").Value = 129
    $ws.Range("O31").Value = 3
    $ws.Range("P31").Value = 0
    # Row 32
    $ws.Range("B32").Value = 22
    $ws.Range("C32").Value = 5
    $ws.Range("D32").Value = 0
    $ws.Range("E32").Value = 1
    $ws.Range("F32").Value = 28
    $ws.Range("G32").Value = 3
    $ws.Range("H32").Value = "<3"
    $ws.Range("I32").Value = 3
    $ws.Range("J32").Value = 4
    $ws.Range("K32").Value = 0
    $ws.Range("L32").Value = 14
    $ws.Range("M32").Value = 8
    $ws.Range("N32").Value = 16
    $ws.Range("O32").Value = 4
    $ws.Range("P32").Value = 0
    # Row 33
    $ws.Range("B33").Value = 217
    $ws.Range("C33").Value = 39
    $ws.Range("D33").Value = 15
    $ws.Range("E33").Value = 10
    $ws.Range("F33").Value = 281
    $ws.Range("G33").Value = 27
    $ws.Range("H33").Value = 6
    $ws.Range("I33").Value = 26
    $ws.Range("J33").Value = 15
    $ws.Range("K33").Value = 0
    $ws.Range("L33").Value = 141
    $ws.Range("M33").Value = 69
    $ws.Range("N33").Value = 199
    $ws.Range("O33").Value = 4
    $ws.Range("P33").Value = 0
    # Row 34
    $ws.Range("B34").Value = 78
    $ws.Range("C34").Value = 22
    $ws.Range("D34").Value = 3
    $ws.Range("E34").Value = 7
    $ws.Range("F34").Value = 110
    $ws.Range("G34").Value = 19
    $ws.Range("H34").Value = "<3"
    $ws.Range("I34").Value = 11
    $ws.Range("J34").Value = 13
    $ws.Range("K34").Value = 0
    $ws.Range("L34").Value = 48
    $ws.Range("M34").Value = 31
    $ws.Range("N34").Value = 70
    $ws.Range("O34").Value = 4
    $ws.Range("P34").Value = 0
    # Row 35
    $ws.Range("B35").Value = 15
    $ws.Range("C35").Value = 4
    $ws.Range("D35").Value = 0
    $ws.Range("E35").Value = 0
    $ws.Range("F35").Value = 19
    $ws.Range("G35").Value = 3
    $ws.Range("H35").Value = "<3"
    $ws.Range("I35").Value = 3
    $ws.Range("J35").Value = "<3"
    $ws.Range("K35").Value = 0
    $ws.Range("L35").Value = 6
    $ws.Range("M35").Value = 8
    $ws.Range("N35").Value = 15
    $ws.Range("O35").Value = 0
    $ws.Range("P35").Value = 0

# --- Update view state: scroll position / active selection ---
$ws.Range("G10").Select()
